# aangepast nav Marcels commentaar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("listProjects")

$nl = [char]10

# Rename the KP ZSS "Spoor" entries: drop the dashes / non-breaking space formatting
$ws.Range("A6").Value  = "KP ZSS Spoor I Zeespiegelstijging en Antarctica"
$ws.Range("A7").Value  = "KP ZSS Spoor II Systeemverkenningen"
$ws.Range("A10").Value = "KP ZSS Spoor V Implementatiestrategie"
$ws.Range("A9").Value  = "KP ZSS Spoor IV Langetermijnopties"
$ws.Range("A8").Value  = "KP ZSS Spoor III Signaleringsmethodiek"

# Wrap long cell texts onto two lines (Alt+Enter equivalent)
$ws.Range("A12").Value = "Programma Beoordelings- en " + $nl + "Ontwerpinstrumentarium (BOI 2023)"
$ws.Range("C3").Value  = "zeespiegel(stijging), " + $nl + "bodemdaling"
$ws.Range("C4").Value  = "zeespiegel(stijging), " + $nl + "bodemdaling"
$ws.Range("C12").Value = "extremen zeespiegel," + $nl + "golven"

# Turn wrap text on for the cells that now contain line breaks, and grow the rows
$ws.Range("C3").WrapText = $true
$ws.Range("C4").WrapText = $true
$ws.Range("A12").WrapText = $true
$ws.Range("C12").WrapText = $true

$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8
$ws.Rows.Item(12).RowHeight = 28.8

# Update the active selection on the sheet
$ws.Range("C14").Select()
